# Antisense_MH_list.xlsx edit:
#  - rename the original "All_microhomology" sheet to "Forward"
#  - duplicate it as a second sheet named "Reverse" containing the
#    reverse-strand antisense microhomology rows
#  - make "Reverse" the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---------------------------------------
$wsFwd = $wb.Worksheets.Item(1)
$wsFwd.Name = "Forward"

# --- Duplicate it to build the "Reverse" sheet ------------------------
# Copy(Before, After) - placing it right after "Forward" preserves tab
# order and carries over all existing formatting (fills/styles) as-is.
$wsFwd.Copy($null, $wsFwd)
$wsRev = $wb.Worksheets.Item(2)
$wsRev.Name = "Reverse"

# Columns B (Celltype), C (Breaks) and D (Type) are identical to the
# Forward sheet already (awt/d5, 2dsb, exon_exon) so only A, E, F, G, H,
# I and J need to be overwritten per row.

$wsRev.Cells.Item(2, 1).Value = "EE1R"
$wsRev.Range("E2:E27").Value = "Reverse"

$rows = @(
    @{ Row=2;  Name="EE1R";  Left=25; Right=223; Pattern="CCT";   Seq="ATGAACTTCACCTCGAAGTTCAT";       Ref="GGAAGTTCACGCCGATGAACTTCACCTCGAAGTTCATCA" }
    @{ Row=3;  Name="EE1R";  Left=25; Right=217; Pattern="CCT";   Seq="ATGAACTTCACCTCGAAGTTCAT";       Ref="GGAAGTTCACGCCGATGAACTTCACCTCGAAGTTCATCA" }
    @{ Row=4;  Name="EE2R";  Left=41; Right=218; Pattern="GCCGTCC"; Seq="AGATGAAGCAGCCGTCCTCGAAGTTCA"; Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCCTCGAAGTTCATCA" }
    @{ Row=5;  Name="EE2R";  Left=41; Right=212; Pattern="GCCGTCC"; Seq="AGATGAAGCAGCCGTCCTCGAAGTTCA"; Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCCTCGAAGTTCATCA" }
    @{ Row=6;  Name="EE3R";  Left=11; Right=218; Pattern="GCCG";  Seq="GGAAGTTCACGCCGTCCTCGAAGT";     Ref="GGAAGTTCACGCCGTCCTCGAAGTTCATCA" }
    @{ Row=7;  Name="EE3R";  Left=11; Right=212; Pattern="GCCG";  Seq="GGAAGTTCACGCCGTCCTCGAAGT";     Ref="GGAAGTTCACGCCGTCCTCGAAGTTCATCA" }
    @{ Row=8;  Name="EE4R";  Left=41; Right=215; Pattern="GCCG";  Seq="AGATGAAGCAGCCGCCGTCCTCGA";     Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=9;  Name="EE4R";  Left=41; Right=209; Pattern="GCCG";  Seq="AGATGAAGCAGCCGCCGTCCTCGA";     Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=10; Name="EE5R";  Left=11; Right=215; Pattern="GCCG";  Seq="GGAAGTTCACGCCGCCGTCCTCGA";     Ref="GGAAGTTCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=11; Name="EE5R";  Left=11; Right=209; Pattern="GCCG";  Seq="GGAAGTTCACGCCGCCGTCCTCGA";     Ref="GGAAGTTCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=12; Name="EE6R";  Left=23; Right=212; Pattern="CAC";   Seq="CGATGAACTTCACGCCGCCGTCC";      Ref="GGAAGTTCACGCCGATGAACTTCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=13; Name="EE6R";  Left=23; Right=206; Pattern="CAC";   Seq="CGATGAACTTCACGCCGCCGTCC";      Ref="GGAAGTTCACGCCGATGAACTTCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=14; Name="EE7R";  Left=41; Right=210; Pattern="GCC";   Seq="AGATGAAGCAGCCACGCCGCCGT";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=15; Name="EE7R";  Left=41; Right=204; Pattern="GCC";   Seq="AGATGAAGCAGCCACGCCGCCGT";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=16; Name="EE8R";  Left=11; Right=210; Pattern="GCC";   Seq="GGAAGTTCACGCCACGCCGCCGT";      Ref="GGAAGTTCACGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=17; Name="EE8R";  Left=11; Right=204; Pattern="GCC";   Seq="GGAAGTTCACGCCACGCCGCCGT";      Ref="GGAAGTTCACGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=18; Name="EE9R";  Left=44; Right=207; Pattern="GTC";   Seq="TGAAGCAGCCGTCGCCACGCCGC";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=19; Name="EE9R";  Left=44; Right=201; Pattern="GTC";   Seq="TGAAGCAGCCGTCGCCACGCCGC";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=20; Name="EE10R"; Left=22; Right=202; Pattern="TCAC";  Seq="CCGATGAACTTCACGGTCGCCACG";     Ref="GGAAGTTCACGCCGATGAACTTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=21; Name="EE10R"; Left=22; Right=196; Pattern="TCAC";  Seq="CCGATGAACTTCACGGTCGCCACG";     Ref="GGAAGTTCACGCCGATGAACTTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=22; Name="EE11R"; Left=44; Right=201; Pattern="GTC";   Seq="TGAAGCAGCCGTCACGGTCGCCA";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=23; Name="EE11R"; Left=44; Right=195; Pattern="GTC";   Seq="TGAAGCAGCCGTCACGGTCGCCA";      Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=24; Name="EE12R"; Left=25; Right=196; Pattern="CCT";   Seq="ATGAACTTCACCTGGGTCACGGT";      Ref="GGAAGTTCACGCCGATGAACTTCACCTGGGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=25; Name="EE12R"; Left=25; Right=190; Pattern="CCT";   Seq="ATGAACTTCACCTGGGTCACGGT";      Ref="GGAAGTTCACGCCGATGAACTTCACCTGGGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=26; Name="EE13R"; Left=44; Right=194; Pattern="GTCC";  Seq="TGAAGCAGCCGTCCTGGGTCACGG";     Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCCTGGGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
    @{ Row=27; Name="EE13R"; Left=44; Right=188; Pattern="GTCC";  Seq="TGAAGCAGCCGTCCTGGGTCACGG";     Ref="GGAAGTTCACGCCGATGAACTTCACCTTGTAGATGAAGCAGCCGTCCTGGGTCACGGTCGCCACGCCGCCGTCCTCGAAGTTCATCA" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $wsRev.Cells.Item($row, 1).Value = $r.Name
    $wsRev.Cells.Item($row, 6).Value = $r.Left
    $wsRev.Cells.Item($row, 7).Value = $r.Right
    $wsRev.Cells.Item($row, 8).Value = $r.Pattern
    $wsRev.Cells.Item($row, 9).Value = $r.Seq
    $wsRev.Cells.Item($row, 10).Value = $r.Ref
}

# --- Selection / active tab -------------------------------------------
$wsRev.Range("D29").Select()
$wsRev.Activate()
